$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "_old" -> "_FV2310" and "_new" -> "_FV2404" ---
$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
for ($i = 0; $i -lt $fv2310Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310Headers[$i]
}

# Column K ("diff") is unchanged.

$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404Headers[$i]
}

# --- 2. Turn the data range into an Excel Table ("Table1") ---
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U71"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split/freeze at row 1) ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null

Write-Output "Applied header rename, table creation, and freeze panes."
